$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = -10.92479999999999
$ws.Range("C10").Value = -12.45380000000001
$ws.Range("C12").Value = -14.44770000000001
$ws.Range("E13").Value = 12.0862
$ws.Range("C18").Value = -14.42700000000001
$ws.Range("C25").Value = -10.82669999999999
